$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers to the Excel parser need to be
# forced to Text format first, so they are stored as strings (matching the source
# data format used throughout column D), not auto-converted to numeric values.
$textCells = @("D5", "D11", "D16", "D18", "D23", "D25", "D26", "D27", "D32", "D39", "D47", "D51")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.989.24'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '1.676.86'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '215.24'
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("E6").Value = '  +1.53%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("D11").Value = '0.0886'
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("D12").Value = '1.912.98'
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").Value = '1.692.41'
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").Value = '65.87'
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").Value = '26.991.14'
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("D18").Value = '237.35'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("E19").Value = '  +5.44%  '
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").Value = '9.20'
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("E24").Value = '  -1.71%  '
$ws.Range("D25").Value = '146.03'
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").Value = '7.21'
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("D27").Value = '16.16'
$ws.Range("E27").Value = '  +1.59%  '
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").Value = '  -0.36%  '
$ws.Range("D32").Value = '3.32'
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").Value = '1.480.63'
$ws.Range("E33").Value = '  +2.10%  '
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("E35").Value = '  +4.91%  '
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("E37").Value = '  +2.22%  '
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("D39").Value = '0.905'
$ws.Range("E39").Value = '  +1.00%  '
$ws.Range("E40").Value = '  -3.44%  '
$ws.Range("E41").Value = '  +1.89%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E43").Value = '  +2.02%  '
$ws.Range("E44").Value = '  +2.11%  '
$ws.Range("D45").Value = '1.819.91'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("D47").Value = '90.51'
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("D48").Value = '0.0₆0106'
$ws.Range("E48").Value = '  +2.43%  '
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("E50").Value = '  +1.75%  '
$ws.Range("D51").Value = '7.75'
$ws.Range("E51").Value = '  +1.21%  '

# Restore the default (Normal) cell style so no stray number-format styling
# is left behind on the cells we temporarily forced to Text.
foreach ($cell in $textCells) {
    $ws.Range($cell).Style = "Normal"
}
